$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 5 workout rows pulled from the Strava sync (rows 185-195).
$newRows = @(
    @{ Row=185; A="Matt";     B=45482; C="Walk";    D=50;  E=2.24; F=174;  G=42;  H=1;  I=0;  J=0;  K=0; L="Sauntering Hippo"; M=5 }
    @{ Row=186; A="Steven";   B=45482; C="Workout"; D=23;  E=0;    F=0;    G=21;  H=2;  I=0;  J=0;  K=0; L="Brave Leopard";    M=5 }
    @{ Row=187; A="Jeremiah"; B=45482; C="Workout"; D=44;  E=0;    F=0;    G=34;  H=10; I=0;  J=0;  K=0; L="Wily Hyena";       M=5 }
    @{ Row=188; A="Eric";     B=45482; C="Run";     D=66;  E=6.06; F=115;  G=0;   H=5;  I=24; J=34; K=0; L="Wily Hyena";       M=5 }
    @{ Row=189; A="Steven";   B=45482; C="Walk";    D=21;  E=0.72; F=138;  G=20;  H=1;  I=0;  J=0;  K=0; L="Brave Leopard";    M=5 }
    @{ Row=190; A="Matt";     B=45483; C="Walk";    D=182; E=5.54; F=912;  G=168; H=4;  I=0;  J=0;  K=0; L="Sauntering Hippo"; M=5 }
    @{ Row=191; A="Jeremiah"; B=45483; C="Workout"; D=45;  E=0;    F=0;    G=40;  H=5;  I=0;  J=0;  K=0; L="Wily Hyena";       M=5 }
    @{ Row=192; A="Steven";   B=45483; C="Workout"; D=31;  E=0;    F=0;    G=18;  H=12; I=1;  J=0;  K=0; L="Brave Leopard";    M=5 }
    @{ Row=193; A="Phil";     B=45483; C="Run";     D=30;  E=2.93; F=233;  G=0;   H=13; I=14; J=1;  K=0; L="Sauntering Hippo"; M=5 }
    @{ Row=194; A="Matt";     B=45483; C="Run";     D=52;  E=4;    F=1079; G=0;   H=10; I=8;  J=19; K=6; L="Sauntering Hippo"; M=5 }
    @{ Row=195; A="Steven";   B=45484; C="Walk";    D=20;  E=1;    F=20;   G=20;  H=0;  I=0;  J=0;  K=0; L="Brave Leopard";    M=5 }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    # Copy the date cell format from the row above so the new date cell
    # reuses the existing short-date style instead of minting a new one.
    $ws.Range("B" + ($r - 1)).Copy() | Out-Null
    $ws.Range("B$r").PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
}
$excel.CutCopyMode = $false

# Scroll the frozen pane down and move the selection to the new first
# empty row below the appended data, matching where the author left off.
$ws.Range("A171").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A196").Select() | Out-Null
